# "adjust property of scene" - tweak a couple of scenes' camera offset /
# rotation values on the Scene.xlsx "Sheet1" table, and leave the sheet
# scrolled/selected where the edit was last made.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 ("Demo1" scene): CamOffestPos / CamOffestRot columns (J/K) now
# match what row 3 ("DemoSummer") already uses.
$ws.Range("J2").Value = "0,8,7"
$ws.Range("K2").Value = "45,180"

# Row 6 ("City" / SelectScene data): CamOffestPos / CamOffestRot columns
# (J/K) get a brand-new camera offset + rotation.
$ws.Range("J6").Value = "0,8,-7"
$ws.Range("K6").Value = "45,0"

# Leave the view scrolled over to column E, with K7 as the active cell
# (matches the view state captured when the workbook was last saved).
$excel.ActiveWindow.ScrollColumn = 5
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("K7").Select()
